$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '67.141.02'
$ws.Range('E2').Value = '  +3.07%  '

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.721.10'
$ws.Range('E3').Value = '  +6.08%  '

# Row 4
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.20%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '420.88'
$ws.Range('E5').Value = '  +0.20%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '131.53'
$ws.Range('E6').Value = '  -0.65%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '3.713.88'
$ws.Range('E7').Value = '  +6.13%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.646'
$ws.Range('E8').Value = '  +0.19%  '

# Row 9
$ws.Range('E9').Value = '  +0.05%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.777'
$ws.Range('E10').Value = '  -1.00%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.186'
$ws.Range('E11').Value = '  +14.81%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0000417'
$ws.Range('E12').Value = '  +61.43%  '

# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '43.24'
$ws.Range('E13').Value = '  -0.67%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '10.67'
$ws.Range('E14').Value = '  +6.64%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '4.287.20'
$ws.Range('E15').Value = '  +5.47%  '

# Row 16
$ws.Range('E16').Value = '  -0.82%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '20.85'
$ws.Range('E17').Value = '  +1.52%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.710.48'
$ws.Range('E18').Value = '  +6.19%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '13.40'
$ws.Range('E19').Value = '  +8.20%  '

# Row 20
$ws.Range('E20').Value = '  +3.68%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '66.995.10'
$ws.Range('E21').Value = '  +3.06%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '448.48'
$ws.Range('E22').Value = '  -3.08%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '16.35'
$ws.Range('E23').Value = '  +22.72%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '90.22'
$ws.Range('E24').Value = '  -0.24%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.16'
$ws.Range('E25').Value = '  -1.65%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '37.80'
$ws.Range('E26').Value = '  +11.38%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.17'
$ws.Range('E27').Value = '  +2.29%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '3.33'
$ws.Range('E28').Value = '  -0.85%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '5.07'
$ws.Range('E29').Value = '  +4.41%  '

# Row 30
$ws.Range('B30').Value = 'Cosmos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '12.77'
$ws.Range('E30').Value = '  +2.36%  '

# Row 31
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.124'
$ws.Range('E31').Value = '  +8.87%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '2.76'
$ws.Range('E32').Value = '  +2.29%  '

# Row 33
$ws.Range('E33').Value = '  -3.01%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.165'
$ws.Range('E34').Value = '  +1.79%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '42.12'
$ws.Range('E35').Value = '  +4.87%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '57.37'
$ws.Range('E36').Value = '  -0.55%  '

# Row 37
$ws.Range('E37').Value = '  +0.11%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0497'
$ws.Range('E38').Value = '  -0.79%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0₃0762'
$ws.Range('E39').Value = '  +8.63%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '3.09'
$ws.Range('E40').Value = '  +31.80%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.149'
$ws.Range('E41').Value = '  +2.84%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '29.79'
$ws.Range('E42').Value = '  +37.55%  '

# Row 43
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.996'
$ws.Range('E43').Value = '  -0.22%  '

# Row 44
$ws.Range('B44').Value = 'LidoDAOToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '3.46'
$ws.Range('E44').Value = '  +4.06%  '

# Row 45
$ws.Range('B45').Value = 'Monero'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '149.01'
$ws.Range('E45').Value = '  +1.85%  '

# Row 46
$ws.Range('B46').Value = 'ARBITRUM'
$ws.Range('C46').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.13'
$ws.Range('E46').Value = '  +4.74%  '

# Row 47
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.17'
$ws.Range('E47').Value = '  +29.68%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.67'
$ws.Range('E48').Value = '  -3.18%  '

# Row 49
$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.91'
$ws.Range('E49').Value = '  -7.03%  '

# Row 50
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '4.39'
$ws.Range('E50').Value = '  -2.89%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.310'
$ws.Range('E51').Value = '  -2.40%  '
